$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated symbol list values scraped on Sat Feb  4 04:04:59 UTC 2023
$ws.Range('D2').Value = "'329.95"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'2.71%"
$ws.Range('E2').Style = 'Normal'
$ws.Range('G2').Value = "'4"
$ws.Range('G2').Style = 'Normal'
$ws.Range('D3').Value = "'41.17"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'3.89%"
$ws.Range('E3').Style = 'Normal'
$ws.Range('G3').Value = "'4"
$ws.Range('G3').Style = 'Normal'
$ws.Range('D4').Value = "'5.702"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'-3.13%"
$ws.Range('E4').Style = 'Normal'
$ws.Range('G4').Value = "'4"
$ws.Range('G4').Style = 'Normal'
$ws.Range('D5').Value = "'0.08166"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'1.92%"
$ws.Range('E5').Style = 'Normal'
$ws.Range('G5').Value = "'4"
$ws.Range('G5').Style = 'Normal'
$ws.Range('D6').Value = "'2.041"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'6.06%"
$ws.Range('E6').Style = 'Normal'
$ws.Range('G6').Value = "'4"
$ws.Range('G6').Style = 'Normal'
$ws.Range('D7').Value = "'8.771"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'1.31%"
$ws.Range('E7').Style = 'Normal'
$ws.Range('G7').Value = "'4"
$ws.Range('G7').Style = 'Normal'
$ws.Range('B8').Value = "'GateToken"
$ws.Range('B8').Style = 'Normal'
$ws.Range('C8').Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range('C8').Style = 'Normal'
$ws.Range('D8').Value = "'4.545"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'-0.50%"
$ws.Range('E8').Style = 'Normal'
$ws.Range('G8').Value = "'4"
$ws.Range('G8').Style = 'Normal'
$ws.Range('B9').Value = "'BTSEToken"
$ws.Range('B9').Style = 'Normal'
$ws.Range('C9').Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range('C9').Style = 'Normal'
$ws.Range('D9').Value = "'2.952"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'0.34%"
$ws.Range('E9').Style = 'Normal'
$ws.Range('G9').Value = "'4"
$ws.Range('G9').Style = 'Normal'
$ws.Range('B10').Value = "'MXToken"
$ws.Range('B10').Style = 'Normal'
$ws.Range('C10').Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range('C10').Style = 'Normal'
$ws.Range('D10').Value = "'0.9243"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'-0.93%"
$ws.Range('E10').Style = 'Normal'
$ws.Range('G10').Value = "'4"
$ws.Range('G10').Style = 'Normal'
$ws.Range('B11').Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range('B11').Style = 'Normal'
$ws.Range('C11').Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range('C11').Style = 'Normal'
$ws.Range('D11').Value = "'0.1249"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'0.92%"
$ws.Range('E11').Style = 'Normal'
$ws.Range('G11').Value = "'4"
$ws.Range('G11').Style = 'Normal'
$ws.Range('B12').Value = "'WazirX"
$ws.Range('B12').Style = 'Normal'
$ws.Range('C12').Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range('C12').Style = 'Normal'
$ws.Range('D12').Value = "'0.1956"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'-0.06%"
$ws.Range('E12').Style = 'Normal'
$ws.Range('G12').Value = "'4"
$ws.Range('G12').Style = 'Normal'
$ws.Range('B13').Value = "'MandalaExchangeToken"
$ws.Range('B13').Style = 'Normal'
$ws.Range('C13').Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range('C13').Style = 'Normal'
$ws.Range('D13').Value = "'0.09305"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'1.16%"
$ws.Range('E13').Style = 'Normal'
$ws.Range('G13').Value = "'4"
$ws.Range('G13').Style = 'Normal'
$ws.Range('B14').Value = "'BitrueCoin"
$ws.Range('B14').Style = 'Normal'
$ws.Range('C14').Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range('C14').Style = 'Normal'
$ws.Range('D14').Value = "'0.03724"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'4.98%"
$ws.Range('E14').Style = 'Normal'
$ws.Range('G14').Value = "'4"
$ws.Range('G14').Style = 'Normal'
$ws.Range('B15').Value = "'BitMartToken"
$ws.Range('B15').Style = 'Normal'
$ws.Range('C15').Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range('C15').Style = 'Normal'
$ws.Range('D15').Value = "'0.1054"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'10.15%"
$ws.Range('E15').Style = 'Normal'
$ws.Range('G15').Value = "'4"
$ws.Range('G15').Style = 'Normal'
$ws.Range('B16').Value = "'BitForexToken"
$ws.Range('B16').Style = 'Normal'
$ws.Range('C16').Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range('C16').Style = 'Normal'
$ws.Range('D16').Value = "'0.001298"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'0.23%"
$ws.Range('E16').Style = 'Normal'
$ws.Range('G16').Value = "'4"
$ws.Range('G16').Style = 'Normal'
$ws.Range('B17').Value = "'TigerCash"
$ws.Range('B17').Style = 'Normal'
$ws.Range('C17').Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range('C17').Style = 'Normal'
$ws.Range('D17').Value = "'0.006136"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'-0.41%"
$ws.Range('E17').Style = 'Normal'
$ws.Range('G17').Value = "'4"
$ws.Range('G17').Style = 'Normal'
$ws.Range('B18').Value = "'LEO"
$ws.Range('B18').Style = 'Normal'
$ws.Range('C18').Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range('C18').Style = 'Normal'
$ws.Range('D18').Value = "'3.434"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'2.47%"
$ws.Range('E18').Style = 'Normal'
$ws.Range('G18').Value = "'4"
$ws.Range('G18').Style = 'Normal'
$ws.Range('G19').Value = "'4"
$ws.Range('G19').Style = 'Normal'
$ws.Range('D20').Value = "'8.327"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'-4.98%"
$ws.Range('E20').Style = 'Normal'
$ws.Range('G20').Value = "'4"
$ws.Range('G20').Style = 'Normal'
$ws.Range('E21').Value = "'-3.34%"
$ws.Range('E21').Style = 'Normal'
$ws.Range('G21').Value = "'4"
$ws.Range('G21').Style = 'Normal'
$ws.Range('G22').Value = "'4"
$ws.Range('G22').Style = 'Normal'
$ws.Range('D23').Value = "'0.04449"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'0.37%"
$ws.Range('E23').Style = 'Normal'
$ws.Range('G23').Value = "'4"
$ws.Range('G23').Style = 'Normal'
$ws.Range('E24').Value = "'0.64%"
$ws.Range('E24').Style = 'Normal'
$ws.Range('G24').Value = "'4"
$ws.Range('G24').Style = 'Normal'
$ws.Range('D25').Value = "'0.004292"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'-2.08%"
$ws.Range('E25').Style = 'Normal'
$ws.Range('G25').Value = "'4"
$ws.Range('G25').Style = 'Normal'
$ws.Range('E26').Value = "'3.69%"
$ws.Range('E26').Style = 'Normal'
$ws.Range('G26').Value = "'4"
$ws.Range('G26').Style = 'Normal'
$ws.Range('G27').Value = "'4"
$ws.Range('G27').Style = 'Normal'
$ws.Range('G28').Value = "'4"
$ws.Range('G28').Style = 'Normal'
$ws.Range('G29').Value = "'4"
$ws.Range('G29').Style = 'Normal'
$ws.Range('G30').Value = "'4"
$ws.Range('G30').Style = 'Normal'
$ws.Range('G31').Value = "'4"
$ws.Range('G31').Style = 'Normal'
$ws.Range('G32').Value = "'4"
$ws.Range('G32').Style = 'Normal'
$ws.Range('G33').Value = "'4"
$ws.Range('G33').Style = 'Normal'
$ws.Range('G34').Value = "'4"
$ws.Range('G34').Style = 'Normal'
$ws.Range('G35').Value = "'4"
$ws.Range('G35').Style = 'Normal'
$ws.Range('G36').Value = "'4"
$ws.Range('G36').Style = 'Normal'
$ws.Range('G37').Value = "'4"
$ws.Range('G37').Style = 'Normal'
$ws.Range('G38').Value = "'4"
$ws.Range('G38').Style = 'Normal'
$ws.Range('D39').Value = "'0.02816"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'16.89%"
$ws.Range('E39').Style = 'Normal'
$ws.Range('G39').Value = "'4"
$ws.Range('G39').Style = 'Normal'
$ws.Range('D40').Value = "'0.05482"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'5.69%"
$ws.Range('E40').Style = 'Normal'
$ws.Range('G40').Value = "'4"
$ws.Range('G40').Style = 'Normal'
$ws.Range('D41').Value = "'0.007673"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'3.47%"
$ws.Range('E41').Style = 'Normal'
$ws.Range('G41').Value = "'4"
$ws.Range('G41').Style = 'Normal'
$ws.Range('D42').Value = "'0.009464"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'0.92%"
$ws.Range('E42').Style = 'Normal'
$ws.Range('G42').Value = "'4"
$ws.Range('G42').Style = 'Normal'
$ws.Range('D43').Value = "'0.1418"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'1.06%"
$ws.Range('E43').Style = 'Normal'
$ws.Range('G43').Value = "'4"
$ws.Range('G43').Style = 'Normal'
$ws.Range('D44').Value = "'0.002124"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'0.17%"
$ws.Range('E44').Style = 'Normal'
$ws.Range('G44').Value = "'4"
$ws.Range('G44').Style = 'Normal'
$ws.Range('D45').Value = "'0.01261"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'12.27%"
$ws.Range('E45').Style = 'Normal'
$ws.Range('G45').Value = "'4"
$ws.Range('G45').Style = 'Normal'
$ws.Range('D46').Value = "'0.00006866"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'1.70%"
$ws.Range('E46').Style = 'Normal'
$ws.Range('G46').Value = "'4"
$ws.Range('G46').Style = 'Normal'
$ws.Range('E47').Value = "'0.17%"
$ws.Range('E47').Style = 'Normal'
$ws.Range('G47').Value = "'4"
$ws.Range('G47').Style = 'Normal'
$ws.Range('E48').Value = "'60.52%"
$ws.Range('E48').Style = 'Normal'
$ws.Range('G48').Value = "'4"
$ws.Range('G48').Style = 'Normal'
$ws.Range('D49').Value = "'0.003529"
$ws.Range('D49').Style = 'Normal'
$ws.Range('G49').Value = "'4"
$ws.Range('G49').Style = 'Normal'
$ws.Range('E50').Value = "'0.17%"
$ws.Range('E50').Style = 'Normal'
$ws.Range('G50').Value = "'4"
$ws.Range('G50').Style = 'Normal'
$ws.Range('E51').Value = "'0.17%"
$ws.Range('E51').Style = 'Normal'
$ws.Range('G51').Value = "'4"
$ws.Range('G51').Style = 'Normal'
